$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9558.272000000001
$ws.Range("I9").Value = 350.2
$ws.Range("J9").Value = 17231.666
$ws.Range("K9").Value = 350.2
$ws.Range("L9").Value = 17231.666
$ws.Range("M9").Value = -181.2
$ws.Range("N9").Value = -17569.666

$ws.Range("H12").Value = 1319.2
$ws.Range("I12").Value = 399
$ws.Range("K12").Value = 399
$ws.Range("M12").Value = -229

$ws.Range("H38").Value = 1959.5454
$ws.Range("I38").Value = 676
$ws.Range("J38").Value = 3499.8
$ws.Range("K38").Value = 2028
$ws.Range("L38").Value = 10499.4
$ws.Range("M38").Value = -1656
$ws.Range("N38").Value = -11243.4

$ws.Range("H39").Value = 214.27272
$ws.Range("I39").Value = 126.333336
$ws.Range("K39").Value = 379.000008
$ws.Range("M39").Value = -83.00000799999998

$ws.Range("H141").Value = 11755.111
$ws.Range("I141").Value = 9959.6
$ws.Range("K141").Value = 29878.8
$ws.Range("M141").Value = -24698.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4202.661
$ws.Range("I32").Value = 3116.6155
$ws.Range("J32").Value = 9850.1
$ws.Range("K32").Value = 3116.6155
$ws.Range("L32").Value = 9850.1
$ws.Range("M32").Value = -2829.6155
$ws.Range("N32").Value = -10424.1

$ws.Range("H61").Value = 3161.7734
$ws.Range("I61").Value = 2413.3264
$ws.Range("K61").Value = 2413.3264
$ws.Range("M61").Value = -2201.3264

$ws.Range("H136").Value = 3161.7734
$ws.Range("I136").Value = 2413.3264
$ws.Range("K136").Value = 7239.9792
$ws.Range("M136").Value = -4689.9792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1899.6863
$ws.Range("I20").Value = 1816.2858
$ws.Range("J20").Value = 2082.125
$ws.Range("K20").Value = 1816.2858
$ws.Range("L20").Value = 2082.125
$ws.Range("M20").Value = -1569.2858
$ws.Range("N20").Value = -2576.125

$ws.Range("H82").Value = 12369.046
$ws.Range("J82").Value = 36313.168
$ws.Range("L82").Value = 36313.168
$ws.Range("N82").Value = -37079.168

$ws.Range("H85").Value = 12369.046
$ws.Range("J85").Value = 36313.168
$ws.Range("L85").Value = 36313.168
$ws.Range("N85").Value = -38965.168

$ws.Range("H94").Value = 35012.117
$ws.Range("I94").Value = 432.73685
$ws.Range("K94").Value = 432.73685
$ws.Range("M94").Value = 18.26315

$ws.Range("H97").Value = 9995.833000000001
$ws.Range("I97").Value = 1050
$ws.Range("J97").Value = 14468.75
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 14468.75
$ws.Range("M97").Value = -59
$ws.Range("N97").Value = -16450.75

$ws.Range("H134").Value = 1697.1774
$ws.Range("I134").Value = 1462.4314
$ws.Range("J134").Value = 2785.5454
$ws.Range("K134").Value = 4387.2942
$ws.Range("L134").Value = 8356.636200000001
$ws.Range("M134").Value = -1852.2942
$ws.Range("N134").Value = -13426.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 380
$ws.Range("I7").Value = 382.75
$ws.Range("K7").Value = 382.75
$ws.Range("M7").Value = -269.75

$ws.Range("H16").Value = 5724.1875
$ws.Range("J16").Value = 9104.25
$ws.Range("L16").Value = 9104.25
$ws.Range("N16").Value = -9678.25

$ws.Range("H36").Value = 12606.5
$ws.Range("I36").Value = 7749.3335
$ws.Range("J36").Value = 16249.375
$ws.Range("K36").Value = 7749.3335
$ws.Range("L36").Value = 16249.375
$ws.Range("M36").Value = -7361.3335
$ws.Range("N36").Value = -17025.375

$ws.Range("H40").Value = 12606.5
$ws.Range("I40").Value = 7749.3335
$ws.Range("J40").Value = 16249.375
$ws.Range("K40").Value = 7749.3335
$ws.Range("L40").Value = 16249.375
$ws.Range("M40").Value = -7589.3335
$ws.Range("N40").Value = -16569.375

$ws.Range("H113").Value = 5724.1875
$ws.Range("J113").Value = 9104.25
$ws.Range("L113").Value = 9104.25
$ws.Range("N113").Value = -13444.25

$ws.Range("H134").Value = 2625.4546
$ws.Range("I134").Value = 2651.628
$ws.Range("K134").Value = 7954.884
$ws.Range("M134").Value = -5419.884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1062.6
$ws.Range("I68").Value = 1009.7778
$ws.Range("J68").Value = 1141.8334
$ws.Range("K68").Value = 3029.3334
$ws.Range("L68").Value = 3425.5002
$ws.Range("M68").Value = -2218.3334
$ws.Range("N68").Value = -5047.5002

$ws.Range("H71").Value = 1062.6
$ws.Range("I71").Value = 1009.7778
$ws.Range("J71").Value = 1141.8334
$ws.Range("K71").Value = 9088.0002
$ws.Range("L71").Value = 10276.5006
$ws.Range("M71").Value = -5032.0002
$ws.Range("N71").Value = -18388.5006

$ws.Range("H122").Value = 149.22223
$ws.Range("J122").Value = 156.125
$ws.Range("L122").Value = 1405.125
$ws.Range("N122").Value = -6305.125

$ws.Range("H137").Value = 1935.48
$ws.Range("I137").Value = 1248.4286
$ws.Range("K137").Value = 3745.2858
$ws.Range("M137").Value = 1354.7142

$ws.Range("H139").Value = 3115.6
$ws.Range("I139").Value = 946.7619
$ws.Range("K139").Value = 2840.2857
$ws.Range("M139").Value = 2299.7143

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6383.6924
$ws.Range("I70").Value = 5844.364
$ws.Range("J70").Value = 6779.2
$ws.Range("K70").Value = 5844.364
$ws.Range("L70").Value = 6779.2
$ws.Range("M70").Value = -5574.364
$ws.Range("N70").Value = -7319.2

$ws.Range("H73").Value = 6383.6924
$ws.Range("I73").Value = 5844.364
$ws.Range("J73").Value = 6779.2
$ws.Range("K73").Value = 5844.364
$ws.Range("L73").Value = 6779.2
$ws.Range("M73").Value = -4908.364
$ws.Range("N73").Value = -8651.200000000001

$ws.Range("H132").Value = 3187.1667
$ws.Range("I132").Value = 3649.889
$ws.Range("J132").Value = 1799
$ws.Range("K132").Value = 10949.667
$ws.Range("L132").Value = 5397
$ws.Range("M132").Value = -8419.667000000001
$ws.Range("N132").Value = -10457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14616.182
$ws.Range("I40").Value = 16249.482
$ws.Range("J40").Value = 2774.75
$ws.Range("K40").Value = 16249.482
$ws.Range("L40").Value = 2774.75
$ws.Range("M40").Value = -16113.482
$ws.Range("N40").Value = -3046.75

$ws.Range("H82").Value = 850.8261
$ws.Range("I82").Value = 690.4
$ws.Range("J82").Value = 1151.625
$ws.Range("K82").Value = 690.4
$ws.Range("L82").Value = 1151.625
$ws.Range("M82").Value = -329.4
$ws.Range("N82").Value = -1873.625

$ws.Range("H85").Value = 850.8261
$ws.Range("I85").Value = 690.4
$ws.Range("J85").Value = 1151.625
$ws.Range("K85").Value = 690.4
$ws.Range("L85").Value = 1151.625
$ws.Range("M85").Value = 557.6
$ws.Range("N85").Value = -3647.625

$ws.Range("H132").Value = 2501.62
$ws.Range("I132").Value = 2150.7659
$ws.Range("J132").Value = 7998.3335
$ws.Range("K132").Value = 6452.297699999999
$ws.Range("L132").Value = 23995.0005
$ws.Range("M132").Value = -3922.297699999999
$ws.Range("N132").Value = -29055.0005

$ws.Range("H136").Value = 20906.936
$ws.Range("I136").Value = 1538.7142
$ws.Range("J136").Value = 183600
$ws.Range("K136").Value = 4616.142599999999
$ws.Range("L136").Value = 550800
$ws.Range("M136").Value = -2066.142599999999
$ws.Range("N136").Value = -555900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 261.1
$ws.Range("I107").Value = 256.77777
$ws.Range("K107").Value = 770.33331
$ws.Range("M107").Value = 1149.66669

$ws.Range("H136").Value = 77101.8
$ws.Range("I136").Value = 89681.27
$ws.Range("J136").Value = 67217.92999999999
$ws.Range("K136").Value = 269043.81
$ws.Range("L136").Value = 201653.79
$ws.Range("M136").Value = -266493.81
$ws.Range("N136").Value = -206753.79
